$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to remain text so values like
# "0.06410" or "26.305.76" are not auto-coerced into numbers and lose
# formatting (trailing zeros, multiple dots, leading/trailing spaces, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.305.76'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.665.22'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  +0.87%  '
$ws.Range("D5").Value = '219.17'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '0.5348'
$ws.Range("E6").Value = '  +1.66%  '
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("D8").Value = '0.2661'
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("D9").Value = '0.06410'
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").Value = '20.61'
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").Value = '0.07820'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '4.563'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = '1.666.05'
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = '1.892.21'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '0.5532'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").Value = '0.0₅8202'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '65.80'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = '1.011'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").Value = '4.694'
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("D20").Value = '194.04'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '10.26'
$ws.Range("E21").Value = '  +1.99%  '
$ws.Range("D22").Value = '6.046'
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").Value = '146.27'
$ws.Range("E24").Value = '  +3.23%  '
$ws.Range("D25").Value = '0.1231'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").Value = '7.200'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").Value = '16.11'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("E28").Value = '  +3.99%  '
$ws.Range("D29").Value = '0.05831'
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("D30").Value = '1.282'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").Value = '3.615'
$ws.Range("E31").Value = '  +2.21%  '
$ws.Range("D32").Value = '3.289'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = '0.9671'
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = '0.5810'
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '0.8686'
$ws.Range("E39").Value = '  +2.79%  '
$ws.Range("D40").Value = '5.871'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '1.054.39'
$ws.Range("E41").Value = '  +3.02%  '
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '104.93'
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("D44").Value = '1.803.34'
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").Value = '57.98'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -4.23%  '
$ws.Range("D47").Value = '1.015'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("E48").Value = '  +1.91%  '
$ws.Range("D49").Value = '8.039'
$ws.Range("E49").Value = '  +2.67%  '
$ws.Range("D50").Value = '0.05166'
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = '1.416'
$ws.Range("E51").Value = '  -3.53%  '

Write-Output "Updated cryptos list"
